# ProductBacklog.xlsx - mark the "Se connecter" (login) and "Se deconnecter"
# (logout) backlog items as done, now that Spring Security has been wired in
# (see commit message: "Add login system with Spring Security (W.I.P)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 - US014 "Se connecter": the SpringSecurity note in column D should
# wrap like the rest of the column, and the "done" flag in column E flips
# from FALSE to TRUE.
$ws.Range("D15").WrapText = $true
$ws.Range("E15").Value = $true

# Row 16 - US015 "Se deconnecter": also flip its "done" flag to TRUE.
$ws.Range("E16").Value = $true

# Scroll the sheet down a bit so row 4 is the first visible row (best effort).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
